$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.061.74"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "'1.895.89"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'312.78"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.5015"
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("D8").Value = "'0.3889"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("D9").Value = "'0.09174"
$ws.Range("E9").Value = "  -4.77%  "
$ws.Range("D10").Value = "'1.127"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").Value = "'41.79"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "'6.384"
$ws.Range("E12").Value = "  -2.70%  "
$ws.Range("D13").Value = "'20.78"
$ws.Range("E13").Value = "  -1.77%  "
$ws.Range("D14").Value = "'1.895.12"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "'7.287"
$ws.Range("E15").Value = "  -3.56%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "'92.35"
$ws.Range("E17").Value = "  -1.59%  "
$ws.Range("E18").Value = "  -2.93%  "
$ws.Range("D19").Value = "'0.06642"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "'17.89"
$ws.Range("E20").Value = "  -1.29%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "'6.213"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("D23").Value = "'28.115.96"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "'2.317"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "'2.116.75"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").Value = "'2.576"
$ws.Range("E27").Value = "  -6.31%  "
$ws.Range("D28").Value = "'20.86"
$ws.Range("E28").Value = "  -2.35%  "
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").Value = "'126.20"
$ws.Range("E30").Value = "  -1.80%  "
$ws.Range("D31").Value = "'1.096"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("D33").Value = "'5.606"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("D34").Value = "'3.618"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").Value = "'9.578"
$ws.Range("E35").Value = "  -2.79%  "
$ws.Range("D36").Value = "'0.06586"
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("D37").Value = "'0.02396"
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("D38").Value = "'0.2199"
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").Value = "'1.224"
$ws.Range("E39").Value = "  -4.29%  "
$ws.Range("D40").Value = "'1.278"
$ws.Range("E40").Value = "  +6.99%  "
$ws.Range("D41").Value = "'0.6478"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("D43").Value = "'11.38"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "'0.6073"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "'13.28"
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").Value = "'1.310"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").Value = "'3.680"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").Value = "'121.86"
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("D51").Value = "'1.183"
$ws.Range("E51").Value = "  -2.35%  "
